$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column I (T7: 21/3/2020)
$ws.Range("I1").Value = "T7: 21/3/2020"

# Column I data values (mirrors column H pattern for this week's update)
$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 2
$ws.Range("I4").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("I7").Value = 9
$ws.Range("I8").Value = 0
$ws.Range("I9").Value = 14
$ws.Range("I10").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("I19").Value = 0

# Extend the shared sum formula across to column I
$ws.Range("I20").Formula = "=SUM(I2:I19)"
